# "Signed Off time sheets"
#
# Fill in the supervisor's name, and the employee's sign-off
# initials + date at the bottom of the weekly timesheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor Name field (row 6, next to "Supervisor Name:" label)
$ws.Range("G6").Value = "Ankita Gangotra"

# Employee signature line (initials) under "Employee Signature"
$ws.Range("A27").Value = "A.G"

# Signature date, same day format as the "Total Hours Reported" date above
$ws.Range("D27").Value = Get-Date -Year 2014 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0

# Match the document's on-screen selection after signing off
$ws.Range("D27:E27").Select() | Out-Null
